$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Cells.Item(28, 4).Value = 44434
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 120
$ws.Cells.Item(28, 11).Value = 20000
$ws.Cells.Item(28, 12).Value = 21000
$ws.Cells.Item(28, 13).Value = 20500
$ws.Cells.Item(28, 16).Value = 1025

# Row 29
$ws.Cells.Item(29, 4).Value = 44449
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 120
$ws.Cells.Item(29, 11).Value = 24000
$ws.Cells.Item(29, 12).Value = 25000
$ws.Cells.Item(29, 13).Value = 24500
$ws.Cells.Item(29, 16).Value = 1225

# Row 30
$ws.Cells.Item(30, 4).Value = 44449
$ws.Cells.Item(30, 9).Value = "Segunda"
$ws.Cells.Item(30, 10).Value = 160
$ws.Cells.Item(30, 11).Value = 22000
$ws.Cells.Item(30, 12).Value = 23000
$ws.Cells.Item(30, 13).Value = 22500
$ws.Cells.Item(30, 16).Value = 1125

# Row 31
$ws.Cells.Item(31, 4).Value = 44358
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 140
$ws.Cells.Item(31, 11).Value = 29000
$ws.Cells.Item(31, 12).Value = 30000
$ws.Cells.Item(31, 13).Value = 29500
$ws.Cells.Item(31, 16).Value = 1475

# Row 32
$ws.Cells.Item(32, 4).Value = 44358
$ws.Cells.Item(32, 9).Value = "Segunda"
$ws.Cells.Item(32, 10).Value = 160
$ws.Cells.Item(32, 11).Value = 27000
$ws.Cells.Item(32, 12).Value = 28000
$ws.Cells.Item(32, 13).Value = 27500
$ws.Cells.Item(32, 16).Value = 1375

# Row 33
$ws.Cells.Item(33, 4).Value = 44298
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 140
$ws.Cells.Item(33, 11).Value = 20000
$ws.Cells.Item(33, 12).Value = 21000
$ws.Cells.Item(33, 13).Value = 20500
$ws.Cells.Item(33, 16).Value = 1025

# Row 34
$ws.Cells.Item(34, 4).Value = 44435
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 240
$ws.Cells.Item(34, 11).Value = 16000
$ws.Cells.Item(34, 12).Value = 21000
$ws.Cells.Item(34, 13).Value = 18500
$ws.Cells.Item(34, 16).Value = 925

# Row 35
$ws.Cells.Item(35, 4).Value = 44435
$ws.Cells.Item(35, 9).Value = "Segunda"
$ws.Cells.Item(35, 10).Value = 120
$ws.Cells.Item(35, 11).Value = 14000
$ws.Cells.Item(35, 12).Value = 15000
$ws.Cells.Item(35, 13).Value = 14500
$ws.Cells.Item(35, 16).Value = 725

# Row 36
$ws.Cells.Item(36, 4).Value = 44442
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 120
$ws.Cells.Item(36, 11).Value = 15000
$ws.Cells.Item(36, 12).Value = 16000
$ws.Cells.Item(36, 13).Value = 15500
$ws.Cells.Item(36, 16).Value = 775

# Row 37
$ws.Cells.Item(37, 4).Value = 44162
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 180
$ws.Cells.Item(37, 11).Value = 24000
$ws.Cells.Item(37, 12).Value = 25000
$ws.Cells.Item(37, 13).Value = 24500
$ws.Cells.Item(37, 16).Value = 1225

# Row 38
$ws.Cells.Item(38, 4).Value = 44162
$ws.Cells.Item(38, 9).Value = "Segunda"
$ws.Cells.Item(38, 10).Value = 200
$ws.Cells.Item(38, 11).Value = 21000
$ws.Cells.Item(38, 12).Value = 22000
$ws.Cells.Item(38, 13).Value = 21500
$ws.Cells.Item(38, 16).Value = 1075

# Row 39
$ws.Cells.Item(39, 4).Value = 44302
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 120
$ws.Cells.Item(39, 11).Value = 28000
$ws.Cells.Item(39, 12).Value = 30000
$ws.Cells.Item(39, 13).Value = 29000
$ws.Cells.Item(39, 16).Value = 1450

# Row 40
$ws.Cells.Item(40, 4).Value = 44302
$ws.Cells.Item(40, 9).Value = "Segunda"
$ws.Cells.Item(40, 10).Value = 120
$ws.Cells.Item(40, 11).Value = 19000
$ws.Cells.Item(40, 12).Value = 20000
$ws.Cells.Item(40, 13).Value = 19500
$ws.Cells.Item(40, 16).Value = 975

# Row 41
$ws.Cells.Item(41, 4).Value = 44308
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 120
$ws.Cells.Item(41, 11).Value = 19000
$ws.Cells.Item(41, 12).Value = 20000
$ws.Cells.Item(41, 13).Value = 19500
$ws.Cells.Item(41, 16).Value = 975

# Row 42
$ws.Cells.Item(42, 4).Value = 44498
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 130
$ws.Cells.Item(42, 11).Value = 37000
$ws.Cells.Item(42, 12).Value = 38000
$ws.Cells.Item(42, 13).Value = 37500
$ws.Cells.Item(42, 16).Value = 1875

# Row 43
$ws.Cells.Item(43, 4).Value = 44396
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 120
$ws.Cells.Item(43, 11).Value = 29000
$ws.Cells.Item(43, 12).Value = 30000
$ws.Cells.Item(43, 13).Value = 29500
$ws.Cells.Item(43, 16).Value = 1475

# Row 44
$ws.Cells.Item(44, 4).Value = 44477
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 130
$ws.Cells.Item(44, 11).Value = 26000
$ws.Cells.Item(44, 12).Value = 27000
$ws.Cells.Item(44, 13).Value = 26500
$ws.Cells.Item(44, 16).Value = 1325

# Row 45
$ws.Cells.Item(45, 4).Value = 44477
$ws.Cells.Item(45, 9).Value = "Segunda"
$ws.Cells.Item(45, 10).Value = 140
$ws.Cells.Item(45, 11).Value = 23000
$ws.Cells.Item(45, 12).Value = 24000
$ws.Cells.Item(45, 13).Value = 23500
$ws.Cells.Item(45, 16).Value = 1175

# Row 46
$ws.Cells.Item(46, 4).Value = 44452
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 120
$ws.Cells.Item(46, 11).Value = 25000
$ws.Cells.Item(46, 12).Value = 26000
$ws.Cells.Item(46, 13).Value = 25500
$ws.Cells.Item(46, 16).Value = 1275

# Row 47
$ws.Cells.Item(47, 4).Value = 44452
$ws.Cells.Item(47, 9).Value = "Segunda"
$ws.Cells.Item(47, 10).Value = 120
$ws.Cells.Item(47, 11).Value = 22000
$ws.Cells.Item(47, 12).Value = 23000
$ws.Cells.Item(47, 13).Value = 22500
$ws.Cells.Item(47, 16).Value = 1125

# Row 48
$ws.Cells.Item(48, 4).Value = 44211
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 120
$ws.Cells.Item(48, 11).Value = 24000
$ws.Cells.Item(48, 12).Value = 25000
$ws.Cells.Item(48, 13).Value = 24500
$ws.Cells.Item(48, 16).Value = 1225

# Row 49
$ws.Cells.Item(49, 4).Value = 44260
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 160
$ws.Cells.Item(49, 11).Value = 24000
$ws.Cells.Item(49, 12).Value = 25000
$ws.Cells.Item(49, 13).Value = 24500
$ws.Cells.Item(49, 16).Value = 1225

# Row 50
$ws.Cells.Item(50, 4).Value = 44446
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 120
$ws.Cells.Item(50, 11).Value = 16000
$ws.Cells.Item(50, 12).Value = 17000
$ws.Cells.Item(50, 13).Value = 16500
$ws.Cells.Item(50, 16).Value = 825

# Row 51
$ws.Cells.Item(51, 4).Value = 44323
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 120
$ws.Cells.Item(51, 11).Value = 21000
$ws.Cells.Item(51, 12).Value = 22000
$ws.Cells.Item(51, 13).Value = 21500
$ws.Cells.Item(51, 16).Value = 1075

# Row 52
$ws.Cells.Item(52, 4).Value = 44323
$ws.Cells.Item(52, 9).Value = "Segunda"
$ws.Cells.Item(52, 10).Value = 120
$ws.Cells.Item(52, 11).Value = 18000
$ws.Cells.Item(52, 12).Value = 19000
$ws.Cells.Item(52, 13).Value = 18500
$ws.Cells.Item(52, 16).Value = 925

# Row 53
$ws.Cells.Item(53, 4).Value = 44295
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 120
$ws.Cells.Item(53, 11).Value = 34000
$ws.Cells.Item(53, 12).Value = 35000
$ws.Cells.Item(53, 13).Value = 34500
$ws.Cells.Item(53, 16).Value = 1725

# Row 54
$ws.Cells.Item(54, 4).Value = 44295
$ws.Cells.Item(54, 9).Value = "Segunda"
$ws.Cells.Item(54, 10).Value = 120
$ws.Cells.Item(54, 11).Value = 31000
$ws.Cells.Item(54, 12).Value = 32000
$ws.Cells.Item(54, 13).Value = 31500
$ws.Cells.Item(54, 16).Value = 1575

# Row 55
$ws.Cells.Item(55, 4).Value = 44372
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 370
$ws.Cells.Item(55, 11).Value = 18000
$ws.Cells.Item(55, 12).Value = 19000
$ws.Cells.Item(55, 13).Value = 18459
$ws.Cells.Item(55, 16).Value = 923

# Row 56
$ws.Cells.Item(56, 4).Value = 44372
$ws.Cells.Item(56, 9).Value = "Segunda"
$ws.Cells.Item(56, 10).Value = 150
$ws.Cells.Item(56, 11).Value = 14000
$ws.Cells.Item(56, 12).Value = 15000
$ws.Cells.Item(56, 13).Value = 14333
$ws.Cells.Item(56, 16).Value = 717

# Row 57
$ws.Cells.Item(57, 4).Value = 44403
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 120
$ws.Cells.Item(57, 11).Value = 29000
$ws.Cells.Item(57, 12).Value = 30000
$ws.Cells.Item(57, 13).Value = 29500
$ws.Cells.Item(57, 16).Value = 1475

# Row 58
$ws.Cells.Item(58, 4).Value = 44169
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 160
$ws.Cells.Item(58, 11).Value = 18000
$ws.Cells.Item(58, 12).Value = 20000
$ws.Cells.Item(58, 13).Value = 19000
$ws.Cells.Item(58, 16).Value = 950

# Row 59
$ws.Cells.Item(59, 4).Value = 44421
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 120
$ws.Cells.Item(59, 11).Value = 23000
$ws.Cells.Item(59, 12).Value = 24000
$ws.Cells.Item(59, 13).Value = 23500
$ws.Cells.Item(59, 16).Value = 1175

# Row 60
$ws.Cells.Item(60, 4).Value = 44421
$ws.Cells.Item(60, 9).Value = "Segunda"
$ws.Cells.Item(60, 10).Value = 120
$ws.Cells.Item(60, 11).Value = 21000
$ws.Cells.Item(60, 12).Value = 22000
$ws.Cells.Item(60, 13).Value = 21500
$ws.Cells.Item(60, 16).Value = 1075

# Row 61
$ws.Cells.Item(61, 4).Value = 44239
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 120
$ws.Cells.Item(61, 11).Value = 25000
$ws.Cells.Item(61, 12).Value = 26000
$ws.Cells.Item(61, 13).Value = 25500
$ws.Cells.Item(61, 16).Value = 1275

# Row 62
$ws.Cells.Item(62, 4).Value = 44239
$ws.Cells.Item(62, 9).Value = "Segunda"
$ws.Cells.Item(62, 10).Value = 120
$ws.Cells.Item(62, 11).Value = 23000
$ws.Cells.Item(62, 12).Value = 24000
$ws.Cells.Item(62, 13).Value = 23500
$ws.Cells.Item(62, 16).Value = 1175

# Row 63
$ws.Cells.Item(63, 4).Value = 44176
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 160
$ws.Cells.Item(63, 11).Value = 17000
$ws.Cells.Item(63, 12).Value = 18000
$ws.Cells.Item(63, 13).Value = 17500
$ws.Cells.Item(63, 16).Value = 875

# Row 64
$ws.Cells.Item(64, 4).Value = 44176
$ws.Cells.Item(64, 9).Value = "Segunda"
$ws.Cells.Item(64, 10).Value = 180
$ws.Cells.Item(64, 11).Value = 15000
$ws.Cells.Item(64, 12).Value = 16000
$ws.Cells.Item(64, 13).Value = 15500
$ws.Cells.Item(64, 16).Value = 775

# Row 65
$ws.Cells.Item(65, 4).Value = 44407
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 140
$ws.Cells.Item(65, 11).Value = 33000
$ws.Cells.Item(65, 12).Value = 34000
$ws.Cells.Item(65, 13).Value = 33500
$ws.Cells.Item(65, 16).Value = 1675

# Row 66
$ws.Cells.Item(66, 4).Value = 44407
$ws.Cells.Item(66, 9).Value = "Segunda"
$ws.Cells.Item(66, 10).Value = 120
$ws.Cells.Item(66, 11).Value = 30000
$ws.Cells.Item(66, 12).Value = 31000
$ws.Cells.Item(66, 13).Value = 30500
$ws.Cells.Item(66, 16).Value = 1525

# Row 67
$ws.Cells.Item(67, 4).Value = 44312
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 120
$ws.Cells.Item(67, 11).Value = 20000
$ws.Cells.Item(67, 12).Value = 21000
$ws.Cells.Item(67, 13).Value = 20500
$ws.Cells.Item(67, 16).Value = 1025

# Row 68
$ws.Cells.Item(68, 4).Value = 44522
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 120
$ws.Cells.Item(68, 11).Value = 28000
$ws.Cells.Item(68, 12).Value = 30000
$ws.Cells.Item(68, 13).Value = 29000
$ws.Cells.Item(68, 16).Value = 1450
